$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# New headers for the Gameweeks import feature
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# New row 2 values: "true" must be stored as text (not boolean), Gameweek as a number
$ws.Range("S2").Value = "'true"
$ws.Range("S2").ClearFormats()
$ws.Range("T2").Value = 3
